# Apply the changes described by the diff to the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
# Only the "想去人数" (want-to-go count) column F is updated for a number of rows.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(5, 6).Value = 1009
$ws1.Cells.Item(7, 6).Value = 606
$ws1.Cells.Item(8, 6).Value = 566
$ws1.Cells.Item(9, 6).Value = 1477
$ws1.Cells.Item(11, 6).Value = 1379
$ws1.Cells.Item(13, 6).Value = 495
$ws1.Cells.Item(14, 6).Value = 1676
$ws1.Cells.Item(15, 6).Value = 1370
$ws1.Cells.Item(16, 6).Value = 813
$ws1.Cells.Item(21, 6).Value = 1148
$ws1.Cells.Item(22, 6).Value = 19
$ws1.Cells.Item(24, 6).Value = 25
$ws1.Cells.Item(25, 6).Value = 3566
$ws1.Cells.Item(26, 6).Value = 706
$ws1.Cells.Item(27, 6).Value = 563
$ws1.Cells.Item(28, 6).Value = 1575

# --- Sheet 2: 演出 (Performance) ---
# The event "广州·今泉爱夏  巡演" (2024-04-24, row 2) is no longer listed; remove that
# row and let the rows below shift up.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Delete()

# --- Sheet 3: 本地生活 (Local Life) ---
# Update the want-to-go count for row 2.
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value = 799

# --- Sheet 4: 全部类型 (All Types) ---
# Same "今泉爱夏" event is removed from the combined listing (row 4 here); shift rows up.
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(4).Delete()
